$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H15").Value = 13282083
$ws.Range("I15").Value = 10966228.95555555
$ws.Range("J15").Value = 2315854.044444444
$ws.Range("R15").Value = 230
$ws.Range("H16").Value = 13282083
$ws.Range("I16").Value = 10966228.95555555
$ws.Range("J16").Value = 2315854.044444444
$ws.Range("R16").Value = 230
$ws.Range("H17").Value = 13282083
$ws.Range("I17").Value = 10966228.95555555
$ws.Range("J17").Value = 2315854.044444444
$ws.Range("R17").Value = 230
$ws.Range("R18").Value = 751
$ws.Range("H22").Value = 29107755.25
$ws.Range("I22").Value = 23173893.70555557
$ws.Range("J22").Value = 5933861.544444445
$ws.Range("R22").Value = 372
$ws.Range("H23").Value = 14533633.25
$ws.Range("I23").Value = 11998288.77305555
$ws.Range("J23").Value = 2535344.476944443
$ws.Range("R23").Value = 244
$ws.Range("H26").Value = 112729676.98
$ws.Range("I26").Value = 78712295.02527781
$ws.Range("J26").Value = 34017381.95472223
$ws.Range("R26").Value = 388
$ws.Range("R27").Value = 1227
$ws.Range("H96").Value = 19559137704.43999
$ws.Range("I96").Value = 19274055783.40999
$ws.Range("J96").Value = 174817966.16
$ws.Range("R96").Value = 125
$ws.Range("R97").Value = 58
$ws.Range("H100").Value = 83359934
$ws.Range("I100").Value = 70855943.90000001
$ws.Range("J100").Value = 11579873.1
$ws.Range("K100").Value = 924117
$ws.Range("R100").Value = 9
$ws.Range("R103").Value = 79
$ws.Range("H129").Value = 443245185.54
$ws.Range("I129").Value = 328629353.63
$ws.Range("K129").Value = 114393331.91
$ws.Range("R129").Value = 1516
$ws.Range("R131").Value = 105
$ws.Range("H145").Value = 10989247331.39
$ws.Range("I145").Value = 4649620517.03
$ws.Range("R145").Value = 445
$ws.Range("R146").Value = 2
$ws.Range("H159").Value = 313884512.1100001
$ws.Range("I159").Value = 278665009.0100001
$ws.Range("R159").Value = 159
$ws.Range("R162").Value = 27
$ws.Range("H185").Value = 4471252422.029998
$ws.Range("I185").Value = 4247689796.819999
$ws.Range("J185").Value = 190123989.0600003
$ws.Range("R185").Value = 1966
$ws.Range("R187").Value = 65
$ws.Range("H246").Value = 4941397505.236057
$ws.Range("I246").Value = 4108174033.263645
$ws.Range("J246").Value = 787849728.2024063
$ws.Range("R246").Value = 6796
$ws.Range("H247").Value = 75717851.32870987
$ws.Range("I247").Value = 62795835.86117981
$ws.Range("J247").Value = 12922015.46752992
$ws.Range("R247").Value = 2384
$ws.Range("H250").Value = 158826831.2423202
$ws.Range("I250").Value = 131073872.2602975
$ws.Range("J250").Value = 27752958.98202248
$ws.Range("R250").Value = 5605
$ws.Range("R251").Value = 1179
$ws.Range("H253").Value = 654310849.2190342
$ws.Range("I253").Value = 538223663.5170459
$ws.Range("J253").Value = 116087185.7019877
$ws.Range("R253").Value = 2281
$ws.Range("H254").Value = 24707099.22052039
$ws.Range("I254").Value = 20614136.99314228
$ws.Range("J254").Value = 4092962.227378105
$ws.Range("R254").Value = 2246
$ws.Range("R255").Value = 51
$ws.Range("H267").Value = 7715751927.713041
$ws.Range("I267").Value = 6248912076.97494
$ws.Range("J267").Value = 1446383610.102084
$ws.Range("R267").Value = 9152
$ws.Range("R269").Value = 35
